# Repull data, push all data, mean calculation
# Update the dSF (column F) values for the affected rows to reflect
# the re-pulled data set.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -1
    3  = 2
    4  = 5
    9  = -7
    10 = -4
    12 = -2
    13 = 4
    14 = -2
    15 = 6
    16 = 3
    17 = -2
    19 = 2
    20 = 4
    21 = 2
    22 = 2
    25 = -4
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
